# ---------------------------------------------------------------------------
# worksheet-3.2.docx -- "Additional Assumptions and some rewording"
#
# 1) Reword "consumer" -> "user" throughout the Assumptions list.
# 2) Assumption 7: "owns" -> "supplies"
# 3) Assumption 8: "is able to" -> "can"; trailing clause reworded to refer
#    to the alarm function.
# 4) New assumption inserted right after Assumption 8: "The user can see
#    the LED Alert Light." (becomes the new Assumption 9; everything after
#    it shifts down by one).
# 5) Old "freezer temperatures" assumption (old #9, now #10) reworded to a
#    concrete "32 oF" threshold (degree mark kept superscript).
# 6) Old "pays their utility bills" assumption (old #11, now #12) reworded
#    to "has working household power".
# 7) "temperature above 150oF" -> "temperatures above 150oF" (old #12, now
#    #13).
# 8) Four new assumptions appended at the end of the list (gravity,
#    radiation, operating a smart phone, operating the smart phone app).
#
# The Assumptions paragraphs are originally Word paragraphs 21-35 (1-15),
# with paragraph 36 a trailing empty placeholder paragraph. All of the
# edits below key off that fixed, known layout.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Assumption 8 (paragraph 28): reworded *before* the blanket
#    "consumer" -> "user" pass below, replacing the whole tail of the
#    sentence (starting at "consumer") in a single shot so none of the
#    original grammar-check proofErr markers around "is able to" are left
#    dangling behind: "consumer is able to hear as to fulfill the purpose
#    of audio speakers." -> "user can hear the alarm function."
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(28)
$p8.Range.Find.Execute("consumer is able to hear as to fulfill the purpose of audio speakers.", `
                        $true, $false, $false, $false, $false, `
                        $true, 1, $false, "user can hear the alarm function.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Global rewording: "consumer" -> "user" (every remaining Assumption
#    that used the word "consumer" gets it swapped for "user" -- 9 more
#    occurrences, all inside the Assumptions list; paragraph count is
#    unaffected).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("consumer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "user", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Assumption 7 (paragraph 27): "... user owns batteries ..." -> "...
#    user supplies batteries ..."
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(27)
$p7.Range.Find.Execute("owns", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "supplies", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Paragraph 29 (still "Assumption 9" at this point): "... will not
#    attach the device to anything below freezer temperatures." -> "...
#    32 oF." (degree mark superscript, matching the other temperature
#    assumption further down).
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(29)
$ptext = $p9.Range.Text
$pStart = $p9.Range.Start
$phrase = "freezer temperatures"
$idx = $ptext.IndexOf($phrase)
$phraseStart = $pStart + $idx
$phraseEnd = $phraseStart + $phrase.Length
$rPhrase = $d.Range($phraseStart, $phraseEnd)
$rPhrase.Text = "32"
$afterPos = $phraseStart + 2
$rAfter = $d.Range($afterPos, $afterPos)
$rAfter.InsertAfter(" oF")

$p9 = $d.Paragraphs(29)
$ptext2 = $p9.Range.Text
$pStart2 = $p9.Range.Start
$oIdx = $ptext2.IndexOf(" oF")
$supStart = $pStart2 + $oIdx
$supEnd = $supStart + 2   # covers the space + "o" -- keep "F" normal script
$rSup = $d.Range($supStart, $supEnd)
$rSup.Font.Superscript = $true

# ---------------------------------------------------------------------------
# 6) Paragraph 31 (still "Assumption 11" at this point): "... user pays
#    their utility bills." -> "... user has working household power."
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(31)
$p11.Range.Find.Execute("pays their utility bills.", $true, $false, `
                         $false, $false, $false, $true, 1, $false, `
                         "has working household power.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Paragraph 32 (still "Assumption 12" at this point): "... temperature
#    above 150oF." -> "... temperatures above 150oF."
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs(32)
$p12.Range.Find.Execute("temperature above", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "temperatures above", 2) | Out-Null

# ---------------------------------------------------------------------------
# Renumber paragraphs 29-35 (old Assumptions 9-15) up by one, from the
# bottom up so an already-renumbered label never collides with a
# not-yet-renumbered one. Each Find is scoped to a single paragraph's own
# Range, so there is no risk of touching a neighboring paragraph.
# ---------------------------------------------------------------------------
$renumberMap = @{29 = "Assumption 9:"; 30 = "Assumption 10:"; `
                  31 = "Assumption 11:"; 32 = "Assumption 12:"; `
                  33 = "Assumption 13:"; 34 = "Assumption 14:"; `
                  35 = "Assumption 15:"}
foreach ($i in 35..29) {
    $oldLabel = $renumberMap[$i]
    $newNum = $i - 29 + 10   # 29->10, 30->11, ... 35->16
    $newLabel = "Assumption " + $newNum + ":"
    $pp = $d.Paragraphs($i)
    $pp.Range.Find.Execute($oldLabel, $true, $false, $false, $false, `
                            $false, $true, 1, $false, $newLabel, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 4) New assumption inserted right after (old) Assumption 8 (paragraph 28):
#    "Assumption 9: The user can see the LED Alert Light."
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(28)
$p8.Range.InsertParagraphAfter()
$pLed = $d.Paragraphs(29)
$pLed.Range.Text = "Assumption 9: The user can see the LED Alert Light."

# ---------------------------------------------------------------------------
# 8) Four brand-new assumptions appended at the very end of the list.
#    Drop the trailing (empty, blue, "_GoBack") placeholder paragraph
#    first, then grow the list from the last real Assumption paragraph
#    (which is now "Assumption 16", one paragraph further down than before
#    because of the LED-light insertion above).
# ---------------------------------------------------------------------------
$pLastReal = $d.Paragraphs(36)
$pTrailing = $d.Paragraphs(37)
if ($pTrailing.Range.Text.Trim() -eq "") {
    $pTrailing.Range.Delete()
}

$newAssumptions = @(
    "Assumption 17: The user will not subject the Tracker to accelerated gravity.",
    "Assumption 18: The user will not expose the Tracker to radiation.",
    "Assumption 19: The user can operate a smart phone.",
    "Assumption 20: The user can operate a smart phone application."
)

$pCursor = $d.Paragraphs(36)
foreach ($text in $newAssumptions) {
    $pCursor.Range.InsertParagraphAfter()
    $pCursor = $pCursor.Next()
    $pCursor.Range.Text = $text
}
